$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must stay text even when it looks like a number
# (matches the original cells' inline-string / text storage so Excel's
# automatic number-detection on Range.Value doesn't silently retype the cell).
function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '65.586.06'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '3.392.47'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue "D5" '559.98'
$ws.Range('E5').Value = '  -0.41%  '
Set-TextValue "D6" '175.38'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('D8').Value = '3.384.78'
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('E9').Value = '  +0.03%  '
Set-TextValue "D10" '0.174'
$ws.Range('E10').Value = '  +3.85%  '
$ws.Range('E11').Value = '  +0.14%  '
Set-TextValue "D12" '53.55'
$ws.Range('E12').Value = '  -2.12%  '
$ws.Range('E13').Value = '  -0.10%  '
Set-TextValue "D14" '9.19'
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('D15').Value = '3.941.33'
$ws.Range('E15').Value = '  -0.06%  '
Set-TextValue "D16" '18.24'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').Value = '3.396.28'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').Value = '65.630.80'
$ws.Range('E19').Value = '  +0.77%  '
Set-TextValue "D20" '11.84'
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('E21').Value = '  +0.16%  '
Set-TextValue "D22" '481.35'
$ws.Range('E22').Value = '  +1.60%  '
Set-TextValue "D23" '4.96'
$ws.Range('E23').Value = '  -0.32%  '
Set-TextValue "D24" '90.08'
$ws.Range('E24').Value = '  +3.73%  '
$ws.Range('E25').Value = '  +3.57%  '
Set-TextValue "D26" '4.09'
$ws.Range('E26').Value = '  -1.40%  '
Set-TextValue "D27" '2.91'
$ws.Range('E27').Value = '  +0.57%  '
Set-TextValue "D28" '10.61'
$ws.Range('E28').Value = '  -2.83%  '
$ws.Range('E29').Value = '  -1.70%  '
Set-TextValue "D30" '31.29'
$ws.Range('E30').Value = '  +1.86%  '
Set-TextValue "D31" '6.56'
$ws.Range('E31').Value = '  -2.17%  '
Set-TextValue "D32" '63.54'
$ws.Range('E32').Value = '  +5.12%  '
Set-TextValue "D33" '11.43'
$ws.Range('E33').Value = '  -1.18%  '
Set-TextValue "D34" '571.89'
$ws.Range('E34').Value = '  -2.17%  '
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E36').Value = '  +0.00%  '
Set-TextValue "D37" '3.64'
$ws.Range('E37').Value = '  +3.70%  '
$ws.Range('E38').Value = '  +0.44%  '
Set-TextValue "D39" '35.81'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue "D40" '0.373'
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0744'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('D42').Value = '3.089.67'
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('E43').Value = '  -2.50%  '
Set-TextValue "D44" '0.0416'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  -3.19%  '
$ws.Range('E47').Value = '  -1.28%  '
$ws.Range('E48').Value = '  +0.03%  '
Set-TextValue "D49" '140.42'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('E50').Value = '  -0.50%  '
Set-TextValue "D51" '8.43'
$ws.Range('E51').Value = '  +0.59%  '
